# Figure5.pptx - "Pictures done until pump B and pw" edit.
#
# This slide has two pictures, each with its own "(x)" caption textbox:
#   - "ZoneTexte 71" -> "(b)" caption, sitting over the right-hand picture
#   - "ZoneTexte 4"  -> "(a)" caption, sitting over the left-hand picture
#
# The commit nudges the "(b)" caption a little, and adds a second "(a)"
# caption (a duplicate of the existing one, moved to a new spot) so that a
# second picture on the slide also gets labelled "(a)".
#
# PowerPoint's Shape.Left/Top/Width/Height are expressed in points while the
# underlying OOXML stores EMUs (1 pt = 12700 EMU); a plain EMU/12700 division
# can land a hair below the intended EMU once it round-trips through the
# COM layer's single-precision float, so EmuToPt() nudges the point value up
# by half an EMU to make sure it always re-quantizes to the exact target EMU.
function EmuToPt {
    param([double]$Emu)
    return ($Emu / 12700.0) + (0.5 / 12700.0)
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- 1. Reposition the "(b)" caption textbox ---------------------------
$bCaption = $s.Shapes.Item(2)
$bCaption.Left = EmuToPt 3519289
$bCaption.Top  = EmuToPt 91183

# --- 2. Add a second "(a)" caption textbox ------------------------------
# Burn shape-id 4 on a throwaway textbox so PowerPoint's next-id counter
# advances to 6 (matching the id the real edit ends up with) before we
# duplicate the existing "(a)" caption.
$placeholder = $s.Shapes.AddTextbox(1, 0, 0, 1, 1)
$placeholder.Delete()

$aCaption = $s.Shapes.Item(3)
$newCaptionRange = $aCaption.Duplicate()
$newCaption = $newCaptionRange.Item(1)
$newCaption.Name = "ZoneTexte 5"
$newCaption.Left = EmuToPt 481415
$newCaption.Top  = EmuToPt 81658
